$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $cell = $ws.Range($cellRef)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue "D2" '70.734.42'
$ws.Range("E2").Value = '  -2.57%  '
Set-TextValue "D3" '3.864.88'
$ws.Range("E3").Value = '  -2.74%  '
Set-TextValue "D4" '0.998'
$ws.Range("E4").Value = '  -0.17%  '
Set-TextValue "D5" '591.57'
$ws.Range("E5").Value = '  +1.10%  '
Set-TextValue "D6" '167.27'
$ws.Range("E6").Value = '  +6.47%  '
Set-TextValue "D7" '0.670'
$ws.Range("E7").Value = '  -1.65%  '
$ws.Range("E8").Value = '  +0.19%  '
Set-TextValue "D9" '0.749'
$ws.Range("E9").Value = '  +0.15%  '
Set-TextValue "D10" '0.175'
$ws.Range("E10").Value = '  +3.32%  '
Set-TextValue "D11" '53.52'
$ws.Range("E11").Value = '  +0.86%  '
Set-TextValue "D12" '0.0000321'
$ws.Range("E12").Value = '  -0.19%  '
Set-TextValue "D13" '11.29'
$ws.Range("E13").Value = '  +4.59%  '
Set-TextValue "D14" '4.471.34'
$ws.Range("E14").Value = '  -3.09%  '
Set-TextValue "D15" '3.845.18'
$ws.Range("E15").Value = '  -3.32%  '
Set-TextValue "D16" '20.68'
$ws.Range("E16").Value = '  +1.53%  '
Set-TextValue "D17" '13.82'
$ws.Range("E17").Value = '  -1.54%  '
$ws.Range("E18").Value = '  -5.69%  '
$ws.Range("E19").Value = '  -2.08%  '
Set-TextValue "D20" '70.641.83'
$ws.Range("E20").Value = '  -2.35%  '
Set-TextValue "D21" '435.60'
$ws.Range("E21").Value = '  +0.83%  '
Set-TextValue "D22" '4.73'
$ws.Range("E22").Value = '  +0.65%  '
Set-TextValue "D23" '94.18'
$ws.Range("E23").Value = '  -1.88%  '
Set-TextValue "D24" '3.28'
$ws.Range("E24").Value = '  -4.23%  '
Set-TextValue "D25" '13.79'
$ws.Range("E25").Value = '  -3.59%  '
Set-TextValue "D26" '4.07'
$ws.Range("E26").Value = '  -7.89%  '
Set-TextValue "D27" '11.08'
$ws.Range("E27").Value = '  +0.42%  '
Set-TextValue "D28" '5.91'
$ws.Range("E28").Value = '  -0.38%  '
Set-TextValue "D29" '10.29'
$ws.Range("E29").Value = '  -3.43%  '
Set-TextValue "D30" '35.02'
$ws.Range("E30").Value = '  -3.98%  '
Set-TextValue "D31" '8.02'
$ws.Range("E31").Value = '  +2.99%  '
Set-TextValue "D32" '13.51'
$ws.Range("E32").Value = '  -0.31%  '
Set-TextValue "D33" '48.98'
$ws.Range("E33").Value = '  +0.96%  '
$ws.Range("B34").Value = 'OKB'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue "D34" '70.05'
$ws.Range("E34").Value = '  +2.00%  '
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue "D35" '0.125'
$ws.Range("E35").Value = '  -4.62%  '
Set-TextValue "D36" '0.0₃0976'
$ws.Range("E36").Value = '  +10.84%  '
Set-TextValue "D37" '623.62'
$ws.Range("E37").Value = '  -8.14%  '
Set-TextValue "D38" '0.420'
$ws.Range("E38").Value = '  -3.60%  '
$ws.Range("E39").Value = '  -0.12%  '
$ws.Range("B40").Value = 'dogwifhat'
$ws.Range("C40").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue "D40" '3.32'
$ws.Range("E40").Value = '  +28.17%  '
Set-TextValue "D41" '0.143'
$ws.Range("E41").Value = '  -1.69%  '
$ws.Range("B42").Value = 'FirstDigitalUSD'
$ws.Range("C42").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue "D42" '0.998'
$ws.Range("E42").Value = '  -0.33%  '
$ws.Range("B43").Value = 'ThetaToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
Set-TextValue "D43" '3.26'
$ws.Range("E43").Value = '  -2.39%  '
Set-TextValue "D44" '0.0468'
$ws.Range("E44").Value = '  -3.53%  '
Set-TextValue "D45" '10.09'
$ws.Range("E45").Value = '  -5.81%  '
Set-TextValue "D46" '2.70'
$ws.Range("E46").Value = '  +2.07%  '
Set-TextValue "D47" '0.144'
$ws.Range("E47").Value = '  -3.55%  '
$ws.Range("E48").Value = '  -3.03%  '
$ws.Range("E49").Value = '  -17.06%  '
Set-TextValue "D50" '2.837.95'
$ws.Range("E50").Value = '  +1.83%  '
Set-TextValue "D51" '0.000271'
$ws.Range("E51").Value = '  -0.31%  '
